# Season record columns: Wins, Losses, Ties
# The sheet previously ended at column AC (A1:AC54). We extend it with three
# new columns (AD, AE, AF) holding the team's season record, repeated for
# every player row, since the whole sheet represents one team/season.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 54

# Copy the formatting of the existing last header cell (AC1) onto the new
# header cells so they keep the bold/centered/bordered header style, then
# set their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record (93 wins, 69 losses, 0 ties) for every data row.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 93
    $ws.Cells.Item($r, 31).Value = 69
    $ws.Cells.Item($r, 32).Value = 0
}
